$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("C1").Value() -eq "Consumption Days") {
        $ws.Range("C1").Value = "Consumption Period"
    }
    if ($ws.Range("D1").Value() -eq "Usage (%)") {
        $ws.Range("D1").Value = "Utilisation (%)"
    }
}
